$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 165.16667
$ws.Range("I28").Value = 165.16667
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 165.16667
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 319.83333
$ws.Range("N28").ClearContents()

$ws.Range("H64").Value = 3839.7778
$ws.Range("I64").Value = 4112.9033
$ws.Range("J64").Value = 3235
$ws.Range("K64").Value = 4112.9033
$ws.Range("L64").Value = 3235
$ws.Range("M64").Value = -3864.9033
$ws.Range("N64").Value = -3731

$ws.Range("H67").Value = 3839.7778
$ws.Range("I67").Value = 4112.9033
$ws.Range("J67").Value = 3235
$ws.Range("K67").Value = 4112.9033
$ws.Range("L67").Value = 3235
$ws.Range("M67").Value = -3254.9033
$ws.Range("N67").Value = -4951

$ws.Range("H74").Value = 3926.4285
$ws.Range("I74").Value = 4297
$ws.Range("K74").Value = 4297
$ws.Range("M74").Value = -3361

$ws.Range("H76").Value = 6396.324
$ws.Range("I76").Value = 7361.9165
$ws.Range("J76").Value = 4613.6924
$ws.Range("K76").Value = 7361.9165
$ws.Range("L76").Value = 4613.6924
$ws.Range("M76").Value = -7046.9165
$ws.Range("N76").Value = -5243.6924

$ws.Range("H77").Value = 3926.4285
$ws.Range("I77").Value = 4297
$ws.Range("K77").Value = 21485
$ws.Range("M77").Value = -16805

$ws.Range("H79").Value = 6396.324
$ws.Range("I79").Value = 7361.9165
$ws.Range("J79").Value = 4613.6924
$ws.Range("K79").Value = 7361.9165
$ws.Range("L79").Value = 4613.6924
$ws.Range("M79").Value = -6269.9165
$ws.Range("N79").Value = -6797.6924

$ws.Range("H98").Value = 1586.4445
$ws.Range("I98").Value = 1453.5
$ws.Range("K98").Value = 1453.5
$ws.Range("M98").Value = 44.5

$ws.Range("H107").Value = 506.1875
$ws.Range("I107").Value = 510.64285
$ws.Range("J107").Value = 475
$ws.Range("K107").Value = 510.64285
$ws.Range("L107").Value = 475
$ws.Range("M107").Value = 1409.35715
$ws.Range("N107").Value = -4315

$ws.Range("H122").Value = 1586.4445
$ws.Range("I122").Value = 1453.5
$ws.Range("K122").Value = 4360.5
$ws.Range("M122").Value = -1910.5

$ws.Range("H132").Value = 2514.9524
$ws.Range("I132").Value = 1989
$ws.Range("K132").Value = 5967
$ws.Range("M132").Value = -3437

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10148.208
$ws.Range("I32").Value = 7230.911
$ws.Range("K32").Value = 7230.911
$ws.Range("M32").Value = -6943.911

$ws.Range("H61").Value = 240941.58
$ws.Range("I61").Value = 2072.4092
$ws.Range("J61").Value = 503697.66
$ws.Range("K61").Value = 2072.4092
$ws.Range("L61").Value = 503697.66
$ws.Range("M61").Value = -1860.4092
$ws.Range("N61").Value = -504121.66

$ws.Range("H74").Value = 1457.6444
$ws.Range("I74").Value = 1117.7858
$ws.Range("J74").Value = 2017.4117
$ws.Range("K74").Value = 1117.7858
$ws.Range("L74").Value = 2017.4117
$ws.Range("M74").Value = -243.7858000000001
$ws.Range("N74").Value = -3765.4117

$ws.Range("H77").Value = 1457.6444
$ws.Range("I77").Value = 1117.7858
$ws.Range("J77").Value = 2017.4117
$ws.Range("K77").Value = 5588.929
$ws.Range("L77").Value = 10087.0585
$ws.Range("M77").Value = -1220.929
$ws.Range("N77").Value = -18823.0585

$ws.Range("H86").Value = 37285
$ws.Range("I86").Value = 37285
$ws.Range("K86").Value = 37285
$ws.Range("M86").Value = -36099

$ws.Range("H89").Value = 37285
$ws.Range("I89").Value = 37285
$ws.Range("K89").Value = 111855
$ws.Range("M89").Value = -105927

$ws.Range("H97").Value = 555.6429000000001
$ws.Range("I97").Value = 362.43478
$ws.Range("J97").Value = 1444.4
$ws.Range("K97").Value = 362.43478
$ws.Range("L97").Value = 1444.4
$ws.Range("M97").Value = 133.56522
$ws.Range("N97").Value = -2436.4

$ws.Range("H110").Value = 1482.5238
$ws.Range("I110").Value = 1256.2667
$ws.Range("K110").Value = 1256.2667
$ws.Range("M110").Value = 788.7333000000001

$ws.Range("H136").Value = 240941.58
$ws.Range("I136").Value = 2072.4092
$ws.Range("J136").Value = 503697.66
$ws.Range("K136").Value = 6217.2276
$ws.Range("L136").Value = 1511092.98
$ws.Range("M136").Value = -3667.2276
$ws.Range("N136").Value = -1516192.98

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 15393.4
$ws.Range("I105").Value = 18501.25
$ws.Range("J105").Value = 2962
$ws.Range("K105").Value = 18501.25
$ws.Range("L105").Value = 2962
$ws.Range("M105").Value = -16754.25
$ws.Range("N105").Value = -6456

$ws.Range("H134").Value = 2140.0789
$ws.Range("I134").Value = 1675.7826
$ws.Range("J134").Value = 2852
$ws.Range("K134").Value = 5027.3478
$ws.Range("L134").Value = 8556
$ws.Range("M134").Value = -2492.3478
$ws.Range("N134").Value = -13626

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1606
$ws.Range("I16").Value = 1242.1
$ws.Range("J16").Value = 2515.75
$ws.Range("K16").Value = 1242.1
$ws.Range("L16").Value = 2515.75
$ws.Range("M16").Value = -955.0999999999999
$ws.Range("N16").Value = -3089.75

$ws.Range("H105").Value = 1538.174
$ws.Range("I105").Value = 1649.9375
$ws.Range("K105").Value = 1649.9375
$ws.Range("M105").Value = 97.0625

$ws.Range("H107").Value = 565.4400000000001
$ws.Range("I107").Value = 218.41176
$ws.Range("K107").Value = 218.41176
$ws.Range("M107").Value = 1701.58824

$ws.Range("H113").Value = 1606
$ws.Range("I113").Value = 1242.1
$ws.Range("J113").Value = 2515.75
$ws.Range("K113").Value = 1242.1
$ws.Range("L113").Value = 2515.75
$ws.Range("M113").Value = 927.9000000000001
$ws.Range("N113").Value = -6855.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1714830.4
$ws.Range("I113").Value = 2778338.8
$ws.Range("K113").Value = 8335016.399999999
$ws.Range("M113").Value = -8332846.399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5246.6177
$ws.Range("I70").Value = 5199.517
$ws.Range("K70").Value = 5199.517
$ws.Range("M70").Value = -4929.517

$ws.Range("H73").Value = 5246.6177
$ws.Range("I73").Value = 5199.517
$ws.Range("K73").Value = 5199.517
$ws.Range("M73").Value = -4263.517

$ws.Range("H113").Value = 1538.7693
$ws.Range("J113").Value = 2163.6365
$ws.Range("L113").Value = 2163.6365
$ws.Range("N113").Value = -6503.636500000001

$ws.Range("H122").Value = 1854736.6
$ws.Range("I122").Value = 7409546.5
$ws.Range("J122").Value = 3133.3333
$ws.Range("K122").Value = 22228639.5
$ws.Range("L122").Value = 9399.999899999999
$ws.Range("M122").Value = -22226189.5
$ws.Range("N122").Value = -14299.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1693.1818
$ws.Range("I16").Value = 862.6
$ws.Range("J16").Value = 9999
$ws.Range("K16").Value = 862.6
$ws.Range("L16").Value = 9999
$ws.Range("M16").Value = -692.6
$ws.Range("N16").Value = -10339

$ws.Range("H40").Value = 2079.0527
$ws.Range("I40").Value = 2187.75
$ws.Range("K40").Value = 2187.75
$ws.Range("M40").Value = -2051.75

$ws.Range("H61").Value = 1459.3334
$ws.Range("I61").Value = 1458.5555
$ws.Range("J61").Value = 1461.6666
$ws.Range("K61").Value = 1458.5555
$ws.Range("L61").Value = 1461.6666
$ws.Range("M61").Value = -1256.5555
$ws.Range("N61").Value = -1865.6666

$ws.Range("H93").Value = 688.7778
$ws.Range("I93").Value = 350
$ws.Range("K93").Value = 350
$ws.Range("M93").Value = 898

$ws.Range("H113").Value = 1459.3334
$ws.Range("I113").Value = 1458.5555
$ws.Range("J113").Value = 1461.6666
$ws.Range("K113").Value = 1458.5555
$ws.Range("L113").Value = 1461.6666
$ws.Range("M113").Value = 711.4445000000001
$ws.Range("N113").Value = -5801.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1281.7307
$ws.Range("I132").Value = 857.9474
$ws.Range("J132").Value = 2432
$ws.Range("K132").Value = 2573.8422
$ws.Range("L132").Value = 7296
$ws.Range("M132").Value = -43.84220000000005
$ws.Range("N132").Value = -12356
